$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.052042066530207
$ws.Range("D2").Value = 1.055268628335639
$ws.Range("E2").Value = 1.048900163478011
$ws.Range("F2").Value = 1.063662236074625
$ws.Range("I2").Value = 1.045481499151992
$ws.Range("J2").Value = 1.057066755259953
$ws.Range("K2").Value = 1.058009552368776
$ws.Range("L2").Value = 1.051658720252838
$ws.Range("M2").Value = 1.066380269394518
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.053202855792058
$ws.Range("D3").Value = 1.056180841177848
$ws.Range("E3").Value = 1.049895937688396
$ws.Range("F3").Value = 1.06477227772745
$ws.Range("I3").Value = 1.045816074696081
$ws.Range("J3").Value = 1.057876841161164
$ws.Range("K3").Value = 1.058734920830761
$ws.Range("L3").Value = 1.052466157880417
$ws.Range("M3").Value = 1.06730462547779
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053953946110288
$ws.Range("D4").Value = 1.056771069751885
$ws.Range("E4").Value = 1.050540534742435
$ws.Range("F4").Value = 1.065490930156557
$ws.Range("I4").Value = 1.046031356457702
$ws.Range("J4").Value = 1.058400449160515
$ws.Range("K4").Value = 1.059203617025322
$ws.Range("L4").Value = 1.052988275643255
$ws.Range("M4").Value = 1.067902533878061
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.054269701245172
$ws.Range("D5").Value = 1.057019193805946
$ws.Range("E5").Value = 1.050811586922102
$ws.Range("F5").Value = 1.065793143309403
$ws.Range("I5").Value = 1.046121571142184
$ws.Range("J5").Value = 1.058620437496409
$ws.Range("K5").Value = 1.059400497801696
$ws.Range("L5").Value = 1.053207690850608
$ws.Range("M5").Value = 1.068153844224087
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.054322717763012
$ws.Range("D6").Value = 1.057060854458923
$ws.Range("E6").Value = 1.050857101511107
$ws.Range("F6").Value = 1.06584389161488
$ws.Range("I6").Value = 1.046136701596776
$ws.Range("J6").Value = 1.058657366515479
$ws.Range("K6").Value = 1.059433545628236
$ws.Range("L6").Value = 1.05324452676515
$ws.Range("M6").Value = 1.068196037395194
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053958165257224
$ws.Range("D7").Value = 1.056774385229058
$ws.Range("E7").Value = 1.050544156304801
$ws.Range("F7").Value = 1.065494967982785
$ws.Range("I7").Value = 1.046032563049137
$ws.Range("J7").Value = 1.058403389191132
$ws.Range("K7").Value = 1.059206248380725
$ws.Range("L7").Value = 1.052991207806143
$ws.Range("M7").Value = 1.067905892096575
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.052434365052745
$ws.Range("D8").Value = 1.055576922084553
$ws.Range("E8").Value = 1.04923663540358
$ws.Range("F8").Value = 1.064037301147379
$ws.Range("I8").Value = 1.045594821353603
$ws.Range("J8").Value = 1.05734064621425
$ws.Range("K8").Value = 1.058254832036047
$ws.Range("L8").Value = 1.051931669894107
$ws.Range("M8").Value = 1.066692703930328
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049749049066064
$ws.Range("D9").Value = 1.053466571125703
$ws.Range("E9").Value = 1.046934631300531
$ws.Range("F9").Value = 1.061471607575361
$ws.Range("I9").Value = 1.044814183154003
$ws.Range("J9").Value = 1.055463559655932
$ws.Range("K9").Value = 1.056573204029831
$ws.Range("L9").Value = 1.050061946698135
$ws.Range("M9").Value = 1.0645532727909
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047958644433357
$ws.Range("D10").Value = 1.052059478447176
$ws.Range("E10").Value = 1.045401297752923
$ws.Range("F10").Value = 1.059763067540757
$ws.Range("I10").Value = 1.044287506706357
$ws.Range("J10").Value = 1.054209183296155
$ws.Range("K10").Value = 1.055448662871111
$ws.Range("L10").Value = 1.048813639485408
$ws.Range("M10").Value = 1.063125857977395
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047183317225147
$ws.Range("D11").Value = 1.051450142685977
$ws.Range("E11").Value = 1.044737659174127
$ws.Range("F11").Value = 1.059023699520276
$ws.Range("I11").Value = 1.0440579645218
$ws.Range("J11").Value = 1.053665309286691
$ws.Range("K11").Value = 1.054960899367457
$ws.Range("L11").Value = 1.048272670169405
$ws.Range("M11").Value = 1.062507497647956
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046895314133785
$ws.Range("D12").Value = 1.051223799529
$ws.Range("E12").Value = 1.044491199536418
$ws.Range("F12").Value = 1.058749130852438
$ws.Range("I12").Value = 1.043972478499233
$ws.Range("J12").Value = 1.05346318126612
$ws.Range("K12").Value = 1.054779597089526
$ws.Range("L12").Value = 1.0480716627855
$ws.Range("M12").Value = 1.062277768311133
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.046957092343661
$ws.Range("D13").Value = 1.051272351305911
$ws.Range("E13").Value = 1.044544063925031
$ws.Range("F13").Value = 1.058808023809906
$ws.Range("I13").Value = 1.043990825675208
$ws.Range("J13").Value = 1.053506543356445
$ws.Range("K13").Value = 1.054818492717485
$ws.Range("L13").Value = 1.048114782606151
$ws.Range("M13").Value = 1.06232704798022
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047159511054151
$ws.Range("D14").Value = 1.051431433270364
$ws.Range("E14").Value = 1.044717285823031
$ws.Range("F14").Value = 1.05900100225588
$ws.Range("I14").Value = 1.044050902788451
$ws.Range("J14").Value = 1.053648603549077
$ws.Range("K14").Value = 1.054945915421626
$ws.Range("L14").Value = 1.048256056212459
$ws.Range("M14").Value = 1.062488509025822
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.047284226261455
$ws.Range("D15").Value = 1.05152944775067
$ws.Range("E15").Value = 1.044824019535209
$ws.Range("F15").Value = 1.059119911307669
$ws.Range("I15").Value = 1.044087888611214
$ws.Range("J15").Value = 1.053736117057877
$ws.Range("K15").Value = 1.05502440814964
$ws.Range("L15").Value = 1.048343090609588
$ws.Range("M15").Value = 1.062587984816467
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048010098739865
$ws.Range("D16").Value = 1.052099916853852
$ws.Range("E16").Value = 1.04544534763908
$ws.Range("F16").Value = 1.059812146191882
$ws.Range("I16").Value = 1.044302709292409
$ws.Range("J16").Value = 1.054245263176582
$ws.Range("K16").Value = 1.055481016600855
$ws.Range("L16").Value = 1.048849532432127
$ws.Range("M16").Value = 1.063166890553527
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048465399518386
$ws.Range("D17").Value = 1.052457741944658
$ws.Range("E17").Value = 1.045835171611753
$ws.Range("F17").Value = 1.060246484421029
$ws.Range("I17").Value = 1.044437062083937
$ws.Range("J17").Value = 1.05456444361427
$ws.Range("K17").Value = 1.055767212492498
$ws.Range("L17").Value = 1.049167090712891
$ws.Range("M17").Value = 1.06352994731451
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04873096208793
$ws.Range("D18").Value = 1.05266644995369
$ws.Range("E18").Value = 1.046062578800714
$ws.Range("F18").Value = 1.060499869169617
$ws.Range("I18").Value = 1.044515284262995
$ws.Range("J18").Value = 1.054750546667457
$ws.Range("K18").Value = 1.055934065689892
$ws.Range("L18").Value = 1.049352274279381
$ws.Range("M18").Value = 1.063741685181041
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.048821510975428
$ws.Range("D19").Value = 1.052737613102604
$ws.Range("E19").Value = 1.046140123788879
$ws.Range("F19").Value = 1.060586274070177
$ws.Range("I19").Value = 1.044541931668638
$ws.Range("J19").Value = 1.054813991233132
$ws.Range("K19").Value = 1.055990944761046
$ws.Range("L19").Value = 1.049415409846939
$ws.Range("M19").Value = 1.063813877762383
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04841655075445
$ws.Range("D20").Value = 1.052419351251152
$ws.Range("E20").Value = 1.045793344145191
$ws.Range("F20").Value = 1.060199879637562
$ws.Range("I20").Value = 1.044422662149165
$ws.Range("J20").Value = 1.054530205767277
$ws.Range("K20").Value = 1.055736514662844
$ws.Range("L20").Value = 1.049133024146375
$ws.Range("M20").Value = 1.063490997557358
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0470999041336
$ws.Range("D21").Value = 1.051384587840983
$ws.Range("E21").Value = 1.044666275045923
$ws.Range("F21").Value = 1.058944173110197
$ws.Range("I21").Value = 1.044033217757394
$ws.Range("J21").Value = 1.053606773377671
$ws.Range("K21").Value = 1.054908396071458
$ws.Range("L21").Value = 1.048214456518829
$ws.Range("M21").Value = 1.062440963958439
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.046272005716383
$ws.Range("D22").Value = 1.050733940333858
$ws.Range("E22").Value = 1.043957903672615
$ws.Range("F22").Value = 1.058155039381531
$ws.Range("I22").Value = 1.043787063419849
$ws.Range("J22").Value = 1.053025543806107
$ws.Range("K22").Value = 1.054387000441574
$ws.Range("L22").Value = 1.04763652696841
$ws.Range("M22").Value = 1.061780518328757
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046710897508291
$ws.Range("D23").Value = 1.051078865829261
$ws.Range("E23").Value = 1.044333400114325
$ws.Range("F23").Value = 1.058573338440847
$ws.Range("I23").Value = 1.043917677363586
$ws.Range("J23").Value = 1.05333372464584
$ws.Range("K23").Value = 1.0546634709922
$ws.Range("L23").Value = 1.047942935483396
$ws.Range("M23").Value = 1.062130656767331
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048438623428382
$ws.Range("D24").Value = 1.052436698369245
$ws.Range("E24").Value = 1.045812244085622
$ws.Range("F24").Value = 1.060220938201198
$ws.Range("I24").Value = 1.044429169303208
$ws.Range("J24").Value = 1.054545676591052
$ws.Range("K24").Value = 1.055750385937734
$ws.Range("L24").Value = 1.049148417493436
$ws.Range("M24").Value = 1.063508597359207
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050443296726024
$ws.Range("D25").Value = 1.054012180401447
$ws.Range("E25").Value = 1.047529517705381
$ws.Range("F25").Value = 1.062134559941986
$ws.Range("I25").Value = 1.045017097837979
$ws.Range("J25").Value = 1.057008552853122
$ws.Range("K25").Value = 1.057876841161164
$ws.Range("L25").Value = 1.050545635040738
$ws.Range("M25").Value = 1.065106563847096
